$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column cells to remain text so values like "1.000" or "0.4341" are not
# reinterpreted by Excel as numbers (matches original inline-string text cells).
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = '27.949.18'
$ws.Range("E2").Value = '  -3.29%  '
$ws.Range("D3").Value = '1.856.89'
$ws.Range("E3").Value = '  -2.79%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").Value = '317.19'
$ws.Range("E5").Value = '  -2.32%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").Value = '0.4341'
$ws.Range("E7").Value = '  -5.37%  '
$ws.Range("D8").Value = '0.3687'
$ws.Range("E8").Value = '  -3.33%  '
$ws.Range("D9").Value = '0.07483'
$ws.Range("E9").Value = '  -3.01%  '
$ws.Range("D10").Value = '0.9380'
$ws.Range("E10").Value = '  -4.24%  '
$ws.Range("D11").Value = '21.30'
$ws.Range("D12").Value = '1.894.98'
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").Value = '6.710'
$ws.Range("E13").Value = '  -3.19%  '
$ws.Range("D14").Value = '5.426'
$ws.Range("E14").Value = '  -4.10%  '
$ws.Range("D15").Value = '0.06851'
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '81.42'
$ws.Range("E17").Value = '  -2.82%  '
$ws.Range("D18").Value = '0.000009012'
$ws.Range("E18").Value = '  -4.69%  '
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("D20").Value = '15.88'
$ws.Range("E20").Value = '  -4.55%  '
$ws.Range("D21").Value = '27.905.29'
$ws.Range("E21").Value = '  -3.47%  '
$ws.Range("D22").Value = '5.091'
$ws.Range("E22").Value = '  -4.27%  '
$ws.Range("D23").Value = '10.97'
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("D24").Value = '2.068.75'
$ws.Range("E24").Value = '  -3.66%  '
$ws.Range("D25").Value = '2.009'
$ws.Range("E25").Value = '  -3.99%  '
$ws.Range("D26").Value = '154.01'
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("D27").Value = '18.33'
$ws.Range("E27").Value = '  -3.73%  '
$ws.Range("D28").Value = '5.391'
$ws.Range("E28").Value = '  -4.66%  '
$ws.Range("D29").Value = '113.52'
$ws.Range("E29").Value = '  -3.60%  '
$ws.Range("D30").Value = '1.726'
$ws.Range("E30").Value = '  -7.49%  '
$ws.Range("D31").Value = '0.08954'
$ws.Range("E31").Value = '  -3.62%  '
$ws.Range("D32").Value = '0.8027'
$ws.Range("E32").Value = '  -7.21%  '
$ws.Range("D33").Value = '4.826'
$ws.Range("E33").Value = '  -4.87%  '
$ws.Range("D34").Value = '2.997'
$ws.Range("E34").Value = '  -3.05%  '
$ws.Range("D35").Value = '1.170'
$ws.Range("E35").Value = '  -6.31%  '
$ws.Range("D36").Value = '1.002'
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("E37").Value = '  -3.85%  '
$ws.Range("D38").Value = '0.05430'
$ws.Range("E38").Value = '  -4.99%  '
$ws.Range("E39").Value = '  -3.80%  '
$ws.Range("D40").Value = '2.912'
$ws.Range("E40").Value = '  +1.56%  '
$ws.Range("D41").Value = '0.5224'
$ws.Range("E41").Value = '  -4.76%  '
$ws.Range("D42").Value = '7.019'
$ws.Range("E42").Value = '  -5.28%  '
$ws.Range("E43").Value = '  -4.33%  '
$ws.Range("D44").Value = '8.761'
$ws.Range("E44").Value = '  -5.79%  '
$ws.Range("D45").Value = '0.06721'
$ws.Range("E45").Value = '  -2.64%  '
$ws.Range("D46").Value = '0.4877'
$ws.Range("E46").Value = '  -5.72%  '
$ws.Range("D47").Value = '10.61'
$ws.Range("E47").Value = '  -5.98%  '
$ws.Range("D48").Value = '106.26'
$ws.Range("E48").Value = '  -3.67%  '
$ws.Range("D49").Value = '1.937'
$ws.Range("E49").Value = '  -7.16%  '
$ws.Range("D50").Value = '1.674'
$ws.Range("E50").Value = '  -5.87%  '
$ws.Range("E51").Value = '  -0.16%  '
